# Update the Fgf9-Fgfr3 LR-pairs sheet with refreshed TPM-derived figures.
# The underlying analysis now collapses the "ECs -> MuSCs" / "ECs -> FAPs" /
# "ECs -> Resolving-Mac" signalling edges into a single ECs-sourced block of
# four target clusters (ECs, FAPs, MuSCs, Resolving-Mac), so the old six-row
# table shrinks to four rows and every numeric column is recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two rows that no longer exist in the refreshed export (delete from
# the bottom up so row numbers of the rows we still need stay stable).
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf9"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.083607
$ws.Range("H2").Value = 3.250821
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.618552666666667
$ws.Range("N2").Value = 13.855658
$ws.Range("O2").Value = 0.78434648953826
$ws.Range("P2").Value = 0.78434648953826
$ws.Range("Q2").Value = 5.004695999468667
$ws.Range("R2").Value = 45.042263995218
$ws.Range("S2").Value = 0.78434648953826
$ws.Range("T2").Value = 0.78434648953826

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf9"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.083607
$ws.Range("H3").Value = 3.250821
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.6792986666666666
$ws.Range("N3").Value = 2.037896
$ws.Range("O3").Value = 0.1153620112191035
$ws.Range("P3").Value = 0.1153620112191036
$ws.Range("Q3").Value = 0.7360927902906665
$ws.Range("R3").Value = 6.624835112616
$ws.Range("S3").Value = 0.1153620112191035
$ws.Range("T3").Value = 0.1153620112191036

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf9"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.083607
$ws.Range("H4").Value = 3.250821
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.478937
$ws.Range("N4").Value = 1.436811
$ws.Range("O4").Value = 0.0813355572127976
$ws.Range("P4").Value = 0.08133555721279762
$ws.Range("Q4").Value = 0.518979485759
$ws.Range("R4").Value = 4.670815371831001
$ws.Range("S4").Value = 0.0813355572127976
$ws.Range("T4").Value = 0.08133555721279762

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fgf9"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.083607
$ws.Range("H5").Value = 3.250821
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1116203333333333
$ws.Range("N5").Value = 0.334861
$ws.Range("O5").Value = 0.01895594202983873
$ws.Range("P5").Value = 0.01895594202983874
$ws.Range("Q5").Value = 0.1209525745423333
$ws.Range("R5").Value = 1.088573170881
$ws.Range("S5").Value = 0.01895594202983873
$ws.Range("T5").Value = 0.01895594202983874
